# Update "想去人数" (want-to-go count) figures to the freshly scraped values.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 606
$ws1.Range("F5").Value = 647
$ws1.Range("F6").Value = 797
$ws1.Range("F15").Value = 302
$ws1.Range("F20").Value = 532
$ws1.Range("F22").Value = 505

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 72
$ws2.Range("F10").Value = 44

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 72
$ws4.Range("F7").Value = 606
$ws4.Range("F9").Value = 647
$ws4.Range("F10").Value = 797
$ws4.Range("F23").Value = 302
$ws4.Range("F27").Value = 44
$ws4.Range("F32").Value = 532
$ws4.Range("F34").Value = 505
